# Update "想去人数" (number of people interested) values in column F
# across the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types)
# sheets, matching the regenerated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 259
$ws.Range("F5").Value = 1807
$ws.Range("F7").Value = 546
$ws.Range("F8").Value = 546
$ws.Range("F9").Value = 5029
$ws.Range("F13").Value = 998
$ws.Range("F14").Value = 338
$ws.Range("F15").Value = 1304
$ws.Range("F17").Value = 1887
$ws.Range("F18").Value = 3010
$ws.Range("F19").Value = 1856
$ws.Range("F22").Value = 174
$ws.Range("F23").Value = 85
$ws.Range("F24").Value = 652
$ws.Range("F25").Value = 950
$ws.Range("F26").Value = 319
$ws.Range("F28").Value = 3323
$ws.Range("F29").Value = 1050
$ws.Range("F30").Value = 2573
$ws.Range("F32").Value = 1615
$ws.Range("F33").Value = 3764
$ws.Range("F34").Value = 101
$ws.Range("F35").Value = 904
$ws.Range("F37").Value = 1175
$ws.Range("F38").Value = 10
$ws.Range("F39").Value = 951
$ws.Range("F40").Value = 1213
$ws.Range("F41").Value = 35
$ws.Range("F42").Value = 919
$ws.Range("F43").Value = 596
$ws.Range("F44").Value = 378
$ws.Range("F45").Value = 378

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = 4
$ws.Range("F16").Value = 14

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 259
$ws.Range("F6").Value = 1807
$ws.Range("F8").Value = 546
$ws.Range("F9").Value = 546
$ws.Range("F10").Value = 5029
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = 338
$ws.Range("F16").Value = 1304
$ws.Range("F17").Value = 3010
$ws.Range("F19").Value = 1856
$ws.Range("F22").Value = 174
$ws.Range("F25").Value = 85
$ws.Range("F26").Value = 950
$ws.Range("F27").Value = 319
$ws.Range("F28").Value = 3323
$ws.Range("F29").Value = 14
$ws.Range("F30").Value = 1050
$ws.Range("F32").Value = 2573
$ws.Range("F33").Value = 1615
$ws.Range("F34").Value = 3764
$ws.Range("F36").Value = 101
$ws.Range("F37").Value = 904
$ws.Range("F38").Value = 1175
$ws.Range("F39").Value = 10
$ws.Range("F40").Value = 951
$ws.Range("F42").Value = 1213
$ws.Range("F43").Value = 35
$ws.Range("F44").Value = 920
$ws.Range("F45").Value = 596
$ws.Range("F46").Value = 378
$ws.Range("F49").Value = 3516
